# Re-run of the NATMI ligand-receptor (Fn1-Itgb6) analysis with updated TPM
# data. The new run only produces results for the non-self sending/target
# cluster pairs, so rows 2-7 get refreshed numbers (and new Target-cluster
# labels) while the former self-pair rows (8-10) are dropped entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2109236666666666
$ws.Range("N2").Value = 0.632771
$ws.Range("O2").Value = 0.8951984155054113
$ws.Range("P2").Value = 0.8951984155054113
$ws.Range("Q2").Value = 1.316186459756
$ws.Range("R2").Value = 11.845678137804
$ws.Range("S2").Value = 0.01550690019962559
$ws.Range("T2").Value = 0.01550690019962559

$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.024693
$ws.Range("N3").Value = 0.07407900000000001
$ws.Range("O3").Value = 0.1048015844945887
$ws.Range("P3").Value = 0.1048015844945887
$ws.Range("Q3").Value = 0.154086986844
$ws.Range("R3").Value = 1.386782881596
$ws.Range("S3").Value = 0.001815405035768175
$ws.Range("T3").Value = 0.001815405035768175

$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 345.566579
$ws.Range("H4").Value = 1036.699737
$ws.Range("I4").Value = 0.9592798330716089
$ws.Range("J4").Value = 0.9592798330716091
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2109236666666666
$ws.Range("N4").Value = 0.632771
$ws.Range("O4").Value = 0.8951984155054113
$ws.Range("P4").Value = 0.8951984155054113
$ws.Range("Q4").Value = 72.88816992013632
$ws.Range("R4").Value = 655.9935292812269
$ws.Range("S4").Value = 0.8587457865919998
$ws.Range("T4").Value = 0.858745786592

$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.024693
$ws.Range("N5").Value = 0.07407900000000001
$ws.Range("O5").Value = 0.1048015844945887
$ws.Range("P5").Value = 0.1048015844945887
$ws.Range("Q5").Value = 8.533075535247001
$ws.Range("R5").Value = 76.797679817223
$ws.Range("S5").Value = 0.1005340464796092
$ws.Range("T5").Value = 0.1005340464796092

$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 8.428738666666666
$ws.Range("H6").Value = 25.286216
$ws.Range("I6").Value = 0.02339786169299727
$ws.Range("J6").Value = 0.02339786169299728
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2109236666666666
$ws.Range("N6").Value = 0.632771
$ws.Range("O6").Value = 0.8951984155054113
$ws.Range("P6").Value = 0.8951984155054113
$ws.Range("Q6").Value = 1.777820464948444
$ws.Range("R6").Value = 16.000384184536
$ws.Range("S6").Value = 0.02094572871378592
$ws.Range("T6").Value = 0.02094572871378592

$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 8.428738666666666
$ws.Range("H7").Value = 25.286216
$ws.Range("I7").Value = 0.02339786169299727
$ws.Range("J7").Value = 0.02339786169299728
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.024693
$ws.Range("N7").Value = 0.07407900000000001
$ws.Range("O7").Value = 0.1048015844945887
$ws.Range("P7").Value = 0.1048015844945887
$ws.Range("Q7").Value = 0.208130843896
$ws.Range("R7").Value = 1.873177595064
$ws.Range("S7").Value = 0.002452132979211354
$ws.Range("T7").Value = 0.002452132979211354

# Remove now-obsolete rows 8-10 (data shrank from 9 pairwise rows to 6)
$ws.Range("A8:T10").EntireRow.Delete()
